$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# This script regenerates the handback-status report: the two source files
# tracked by the workbook have been re-handed-off/handed-back under new
# GUIDs, and the corresponding handoff/handback timestamps have moved on.
#   4e02d25f-208e-478f-8b11-99c7ff49e0f5  ->  e847ce95-6c36-4be5-a0aa-06edb57db761
#   fe0e1dab-8ca0-4760-bb1e-5f3ccc084aa5  ->  ffffda95e9fc-5d55-42c8-99b9-06d152471643
# and the handoff/handback xlf file hash
#   6d900131576a7630fa91e1b65db76e878e57c41c / 48da20e1fad582794a40babfa9f822ac02be869b
#   -> 850871e52c5bd0d1165539a49096239432a2d5db (now shared by both rows)
# ---------------------------------------------------------------------------

function Set-CellAndHyperlinkText {
    # NOTE: this runtime's PowerShell-style function calls only reliably
    # bind *positional* arguments, so avoid named parameters (-ws/-map)
    # both in the signature usage below and at call sites.
    param($ws, $map)

    # Update the visible text of any hyperlinks anchored on the cells we care
    # about (this rewrites the <hyperlink display="..."/> attribute without
    # touching the underlying Address/relationship or the cell's style).
    foreach ($hl in $ws.Hyperlinks) {
        $addr = $hl.Range.Address($false, $false)
        if ($map.ContainsKey($addr)) {
            $hl.TextToDisplay = $map[$addr]
        }
    }

    # Update the actual cell text/value to match.
    foreach ($addr in $map.Keys) {
        $ws.Range($addr).Value = $map[$addr]
    }
}

# ---------------------------------------------------------------------------
# Sheet "Overview"
# ---------------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")

Set-CellAndHyperlinkText $wsOverview @{
    "A2" = "e847ce95-6c36-4be5-a0aa-06edb57db761.md"
    "A3" = "ffffda95e9fc-5d55-42c8-99b9-06d152471643.md"
}

# ---------------------------------------------------------------------------
# Sheet "zh-cn"
# ---------------------------------------------------------------------------
$wsZhCn = $wb.Worksheets.Item("zh-cn")

Set-CellAndHyperlinkText $wsZhCn @{
    "A2" = "e847ce95-6c36-4be5-a0aa-06edb57db761.md"
    "F2" = "e847ce95-6c36-4be5-a0aa-06edb57db761.md"
    "D2" = "e847ce95-6c36-4be5-a0aa-06edb57db761.850871e52c5bd0d1165539a49096239432a2d5db.zh-cn.xlf"
    "G2" = "e847ce95-6c36-4be5-a0aa-06edb57db761.850871e52c5bd0d1165539a49096239432a2d5db.zh-cn.xlf"

    "A3" = "ffffda95e9fc-5d55-42c8-99b9-06d152471643.md"
    "F3" = "ffffda95e9fc-5d55-42c8-99b9-06d152471643.md"
    "D3" = "e847ce95-6c36-4be5-a0aa-06edb57db761.850871e52c5bd0d1165539a49096239432a2d5db.zh-cn.xlf"
    "G3" = "e847ce95-6c36-4be5-a0aa-06edb57db761.850871e52c5bd0d1165539a49096239432a2d5db.zh-cn.xlf"
}

# Correspond Handoff/Handback datetimes (plain text cells, no hyperlink)
$wsZhCn.Range("E2").Value = "2016-03-18 16:54:34"
$wsZhCn.Range("H2").Value = "2016-03-18 16:54:51"
$wsZhCn.Range("E3").Value = "2016-03-18 16:54:34"
$wsZhCn.Range("H3").Value = "2016-03-18 16:54:51"

# ---------------------------------------------------------------------------
# Sheet "de-de"
# ---------------------------------------------------------------------------
$wsDeDe = $wb.Worksheets.Item("de-de")

Set-CellAndHyperlinkText $wsDeDe @{
    "A2" = "e847ce95-6c36-4be5-a0aa-06edb57db761.md"
    "F2" = "e847ce95-6c36-4be5-a0aa-06edb57db761.md"
    "D2" = "e847ce95-6c36-4be5-a0aa-06edb57db761.850871e52c5bd0d1165539a49096239432a2d5db.de-de.xlf"
    "G2" = "e847ce95-6c36-4be5-a0aa-06edb57db761.850871e52c5bd0d1165539a49096239432a2d5db.de-de.xlf"

    "A3" = "ffffda95e9fc-5d55-42c8-99b9-06d152471643.md"
    "F3" = "ffffda95e9fc-5d55-42c8-99b9-06d152471643.md"
    "D3" = "e847ce95-6c36-4be5-a0aa-06edb57db761.850871e52c5bd0d1165539a49096239432a2d5db.de-de.xlf"
    "G3" = "e847ce95-6c36-4be5-a0aa-06edb57db761.850871e52c5bd0d1165539a49096239432a2d5db.de-de.xlf"
}

# Correspond Handoff/Handback datetimes (plain text cells, no hyperlink)
$wsDeDe.Range("E2").Value = "2016-03-18 16:54:37"
$wsDeDe.Range("H2").Value = "2016-03-18 16:54:56"
$wsDeDe.Range("E3").Value = "2016-03-18 16:54:37"
$wsDeDe.Range("H3").Value = "2016-03-18 16:54:56"

Write-Output "Handback status report regenerated."
